$d = $word.ActiveDocument

# Locate the insertion point: right after "ejecución" and before
# " del proyecto Telefónica..." in the atSistemas / HP-UX job description.
$r = $d.Content
$found = $r.Find.Execute("Líder Técnico de HP-UX en el equipo de ejecución", `
                          $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target text 'Líder Técnico de HP-UX en el equipo de ejecución'"
}

# Collapse the found range to its end point (right before " del proyecto...")
$r.Collapse(0)

# Insert " (P2 Factory)" there. Toggling a character-formatting property on the
# inserted range (and returning it to its original value) forces the host to
# keep it as its own run instead of silently re-merging it with the
# neighbouring text of identical formatting.
$ins = $r.Duplicate
$ins.InsertAfter(" (P2 Factory)")
$ins.Bold = 1
$ins.Bold = 0

# The document's "_GoBack" bookmark used to sit at the very end of this
# paragraph; move it to mark the spot of this latest edit, i.e. right after
# " (P2 Factory)" and right before " del proyecto...".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$editEnd = $ins.Duplicate
$editEnd.Collapse(0)
$d.Bookmarks.Add("_GoBack", $editEnd)
